$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns A-D (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$data = @(
    @(1, 1, 5, 5),
    @(4, 1, 10, 10),
    @(6, 1, 15, 15),
    @(8, 1, 20, 20),
    @(9, 2, 5, 5),
    @(10, 2, 10, 10),
    @(3, 3, 5, 5),
    @(5, 3, 10, 10),
    @(2, 4, 5, 6),
    @(7, 4, 11, 11)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $row++
}
